$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 137.46153
$ws.Range("I33").Value = 167.5
$ws.Range("K33").Value = 167.5
$ws.Range("M33").Value = 61.5
$ws.Range("H92").Value = 3949.75
$ws.Range("I92").Value = 3949.75
$ws.Range("K92").Value = 3949.75
$ws.Range("M92").Value = -2701.75
$ws.Range("H112").Value = 4392.1113
$ws.Range("J112").Value = 4392.1113
$ws.Range("L112").Value = 13176.3339
$ws.Range("N112").Value = -15392.3339

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1490.75
$ws.Range("I2").Value = 1754.3334
$ws.Range("J2").Value = 700
$ws.Range("K2").Value = 1754.3334
$ws.Range("L2").Value = 700
$ws.Range("M2").Value = -1641.3334
$ws.Range("N2").Value = -926
$ws.Range("H45").Value = 7247.5
$ws.Range("I45").Value = 8996.666999999999
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 8996.666999999999
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -8619.666999999999
$ws.Range("N45").Value = -2754
$ws.Range("H88").Value = 2996.6667
$ws.Range("I88").Value = 3506
$ws.Range("J88").Value = 2742
$ws.Range("K88").Value = 3506
$ws.Range("L88").Value = 2742
$ws.Range("M88").Value = -3100
$ws.Range("N88").Value = -3554
$ws.Range("H91").Value = 2996.6667
$ws.Range("I91").Value = 3506
$ws.Range("J91").Value = 2742
$ws.Range("K91").Value = 3506
$ws.Range("L91").Value = 2742
$ws.Range("M91").Value = -2102
$ws.Range("N91").Value = -5550
$ws.Range("H92").Value = 25000
$ws.Range("J92").Value = 25000
$ws.Range("L92").Value = 25000
$ws.Range("N92").Value = -29992
$ws.Range("H110").Value = 3084323.8
$ws.Range("I110").Value = 4111809.5
$ws.Range("K110").Value = 4111809.5
$ws.Range("M110").Value = -4109764.5
$ws.Range("H116").Value = 1490.75
$ws.Range("I116").Value = 1754.3334
$ws.Range("J116").Value = 700
$ws.Range("K116").Value = 1754.3334
$ws.Range("L116").Value = 700
$ws.Range("M116").Value = 539.6666
$ws.Range("N116").Value = -5288
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1490.75
$ws.Range("I3").Value = 1754.3334
$ws.Range("J3").Value = 700
$ws.Range("K3").Value = 1754.3334
$ws.Range("L3").Value = 700
$ws.Range("M3").Value = -1640.3334
$ws.Range("N3").Value = -928
$ws.Range("H94").Value = 634.4286
$ws.Range("I94").Value = 657
$ws.Range("K94").Value = 657
$ws.Range("M94").Value = -206
$ws.Range("H134").Value = 1416.6666
$ws.Range("I134").Value = 1500
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 4500
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -1965
$ws.Range("N134").Value = -8070

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H113").Value = 1648.5
$ws.Range("I113").Value = 1647.5
$ws.Range("J113").Value = 1649.5
$ws.Range("K113").Value = 4942.5
$ws.Range("L113").Value = 4948.5
$ws.Range("M113").Value = -2772.5
$ws.Range("N113").Value = -9288.5
$ws.Range("H115").Value = 2000
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5528
$ws.Range("I122").Value = 5528
$ws.Range("K122").Value = 16584
$ws.Range("M122").Value = -14134

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9902.4
$ws.Range("I7").Value = 12348.667
$ws.Range("K7").Value = 12348.667
$ws.Range("M7").Value = -12236.667
$ws.Range("H22").Value = 999.6667
$ws.Range("I22").Value = 999.6667
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 999.6667
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -704.6667
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 999.6667
$ws.Range("I27").Value = 999.6667
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 999.6667
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -892.6667
$ws.Range("N27").ClearContents()
$ws.Range("H40").Value = 1300.4286
$ws.Range("I40").Value = 600.5
$ws.Range("J40").Value = 5500
$ws.Range("K40").Value = 600.5
$ws.Range("L40").Value = 5500
$ws.Range("M40").Value = -464.5
$ws.Range("N40").Value = -5772
$ws.Range("H55").Value = 1112.7693
$ws.Range("J55").Value = 1309.4445
$ws.Range("L55").Value = 1309.4445
$ws.Range("N55").Value = -1655.4445
$ws.Range("H61").Value = 936.8570999999999
$ws.Range("I61").Value = 866.6667
$ws.Range("J61").Value = 989.5
$ws.Range("K61").Value = 866.6667
$ws.Range("L61").Value = 989.5
$ws.Range("M61").Value = -664.6667
$ws.Range("N61").Value = -1393.5
$ws.Range("H82").Value = 3487.5
$ws.Range("I82").Value = 3476
$ws.Range("K82").Value = 3476
$ws.Range("M82").Value = -3115
$ws.Range("H85").Value = 3487.5
$ws.Range("I85").Value = 3476
$ws.Range("K85").Value = 3476
$ws.Range("M85").Value = -2228
$ws.Range("H113").Value = 936.8570999999999
$ws.Range("I113").Value = 866.6667
$ws.Range("J113").Value = 989.5
$ws.Range("K113").Value = 866.6667
$ws.Range("L113").Value = 989.5
$ws.Range("M113").Value = 1303.3333
$ws.Range("N113").Value = -5329.5
$ws.Range("H122").Value = 7477.3125
$ws.Range("I122").Value = 4829.6
$ws.Range("J122").Value = 8680.817999999999
$ws.Range("K122").Value = 14488.8
$ws.Range("L122").Value = 26042.454
$ws.Range("M122").Value = -12038.8
$ws.Range("N122").Value = -30942.454
$ws.Range("H126").Value = 9902.4
$ws.Range("I126").Value = 12348.667
$ws.Range("K126").Value = 37046.001
$ws.Range("M126").Value = -34576.001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1429781.2
$ws.Range("I81").Value = 1516.75
$ws.Range("K81").Value = 3033.5
$ws.Range("M81").Value = -1972.5
$ws.Range("H84").Value = 1429781.2
$ws.Range("I84").Value = 1516.75
$ws.Range("K84").Value = 15167.5
$ws.Range("M84").Value = -9863.5
$ws.Range("H126").Value = 3492.4285
$ws.Range("I126").Value = 2949.3333
$ws.Range("K126").Value = 8847.999899999999
$ws.Range("M126").Value = -6377.999899999999
